# Refactoring core library to handle revised single register template
#
# - "testreg" sheet is renamed to "testreg1"
# - "registerinfo" sheet gets two new rows of content (source / maintainer)
# - "testreg1" sheet gets two new columns of content (source / broader) populated
#   for its data rows, plus new reference-data rows
# - the "testreg1" sheet becomes the active tab/selected sheet instead of "registerinfo"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("registerinfo")
$ws2 = $wb.Worksheets.Item("testreg")

# --- Add new register metadata to the "registerinfo" sheet ---------------
$ws1.Range("B6").Value = "Github repo at  https://github.com/CSIRO-LW-LD/ldrpyutils"
$ws1.Range("B7").Value = "Jonathan Yu"

# --- Add new "source"/"broader" data to the "testreg" data sheet ---------
$ws2.Range("F2").Value = "Green book"
$ws2.Range("F3").Value = "Red book"
$ws2.Range("G3").Value = 1
$ws2.Range("F4").Value = "Yellow book"
$ws2.Range("G4").Value = 1

# --- Update the "id" value on registerinfo to match the renamed sheet ----
# (done after the other new strings so the shared-string table ordering
# matches the updated workbook)
$ws1.Range("B2").Value = "testreg1"

# --- Rename the data sheet itself -----------------------------------------
$ws2.Name = "testreg1"

# --- Update selections and active sheet/tab -------------------------------
[void]$ws1.Range("B2").Select()
[void]$ws2.Range("G5").Select()
[void]$ws2.Activate()
